$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 26031
$ws.Range("E2").Value = 479
$ws.Range("F2").Value = 479
$ws.Range("G2").Value = 902
$ws.Range("H2").Value = 759
$ws.Range("I2").Value = 730
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 13208
$ws.Range("L2").Value = 7133
$ws.Range("M2").Value = 6075
$ws.Range("N2").Value = 5560
$ws.Range("O2").Value = 515
$ws.Range("P2").Value = 174
$ws.Range("Q2").Value = 208
$ws.Range("R2").Value = -369
$ws.Range("S2").Value = 35
$ws.Range("T2").Value = 624
$ws.Range("U2").Value = -416
$ws.Range("V2").Value = 862
$ws.Range("W2").Value = 1.84
$ws.Range("X2").Value = 2.92
$ws.Range("Y2").Value = 13.86
$ws.Range("Z2").Value = 6.01
$ws.Range("AA2").Value = 117.43
$ws.Range("AB2").Value = 3007.97
$ws.Range("AC2").Value = 20938
$ws.Range("AD2").Value = 4.92
$ws.Range("AE2").Value = 159574
$ws.Range("AF2").Value = 0.65
$ws.Range("AG2").Value = 1250
$ws.Range("AH2").Value = 1.21
$ws.Range("AI2").Value = 5.97
$ws.Range("AJ2").Value = 3484800

# Row 3
$ws.Range("D3").Value = 15154
$ws.Range("E3").Value = 219
$ws.Range("F3").Value = -133
$ws.Range("G3").Value = 312
$ws.Range("H3").Value = -213
$ws.Range("I3").Value = -155
$ws.Range("J3").Value = -57
$ws.Range("K3").Value = 11678
$ws.Range("L3").Value = 5913
$ws.Range("M3").Value = 5765
$ws.Range("N3").Value = 5345
$ws.Range("O3").Value = 420
$ws.Range("P3").Value = 174
$ws.Range("Q3").Value = -503
$ws.Range("R3").Value = -1284
$ws.Range("S3").Value = 476
$ws.Range("T3").Value = 735
$ws.Range("U3").Value = -1238
$ws.Range("V3").Value = 1483
$ws.Range("W3").Value = 1.44
$ws.Range("X3").Value = -1.4
$ws.Range("Y3").Value = -2.85
$ws.Range("Z3").Value = -1.71
$ws.Range("AA3").Value = 102.58
$ws.Range("AB3").Value = 2891.26
$ws.Range("AC3").Value = -4455
$ws.Range("AD3").Value = -18.58
$ws.Range("AE3").Value = 153421
$ws.Range("AF3").Value = 0.54
$ws.Range("AG3").Value = 1250
$ws.Range("AH3").Value = 1.51
$ws.Range("AI3").Value = -28.05
$ws.Range("AJ3").Value = 3484800

# Row 4
$ws.Range("D4").Value = 1599
$ws.Range("E4").Value = -1
$ws.Range("F4").Value = 287
$ws.Range("G4").Value = -26
$ws.Range("H4").Value = 191
$ws.Range("I4").Value = 180
$ws.Range("J4").Value = 11
$ws.Range("K4").Value = 11261
$ws.Range("L4").Value = 5382
$ws.Range("M4").Value = 5879
$ws.Range("N4").Value = 5437
$ws.Range("O4").Value = 443
$ws.Range("P4").Value = 174
$ws.Range("Q4").Value = 269
$ws.Range("R4").Value = -255
$ws.Range("S4").Value = -57
$ws.Range("T4").Value = 334
$ws.Range("U4").Value = -66
$ws.Range("V4").Value = 1477
$ws.Range("W4").Value = -0.04
$ws.Range("X4").Value = 11.96
$ws.Range("Y4").Value = 3.34
$ws.Range("Z4").Value = 1.67
$ws.Range("AA4").Value = 91.54
$ws.Range("AB4").Value = 2976.47
$ws.Range("AC4").Value = 5161
$ws.Range("AD4").Value = 13.79
$ws.Range("AE4").Value = 163748
$ws.Range("AF4").Value = 0.43
$ws.Range("AG4").Value = 1250
$ws.Range("AH4").Value = 1.76
$ws.Range("AI4").Value = 23.08
$ws.Range("AJ4").Value = 3484800

# Row 5
$ws.Range("D5").Value = 1951
$ws.Range("E5").Value = 364
$ws.Range("F5").Value = 364
$ws.Range("G5").Value = 479
$ws.Range("H5").Value = 519
$ws.Range("I5").Value = 492
$ws.Range("J5").Value = 27
$ws.Range("K5").Value = 5993
$ws.Range("L5").Value = 1989
$ws.Range("M5").Value = 4004
$ws.Range("N5").Value = 3529
$ws.Range("O5").Value = 475
$ws.Range("P5").Value = 118
$ws.Range("Q5").Value = 602
$ws.Range("R5").Value = -202
$ws.Range("S5").Value = -63
$ws.Range("T5").Value = 220
$ws.Range("U5").Value = 382
$ws.Range("V5").Value = 1309
$ws.Range("W5").Value = 18.65
$ws.Range("X5").Value = 26.62
$ws.Range("Y5").Value = 10.98
$ws.Range("Z5").Value = 6.02
$ws.Range("AA5").Value = 49.67
$ws.Range("AB5").Value = 5112.84
$ws.Range("AC5").Value = 19741
$ws.Range("AD5").Value = 2.17
$ws.Range("AE5").Value = 155248
$ws.Range("AF5").Value = 0.28
$ws.Range("AG5").Value = 500
$ws.Range("AH5").Value = 1.17
$ws.Range("AI5").Value = 2.31
$ws.Range("AJ5").Value = 2365023

# Row 6
$ws.Range("D6").Value = 2286
$ws.Range("E6").Value = 392
$ws.Range("F6").Value = 392
$ws.Range("G6").Value = 307
$ws.Range("H6").Value = 298
$ws.Range("I6").Value = 273
$ws.Range("K6").Value = 5920
$ws.Range("L6").Value = 1941
$ws.Range("M6").Value = 3979
$ws.Range("N6").Value = 3895
$ws.Range("P6").Value = 118
$ws.Range("Q6").Value = 125
$ws.Range("R6").Value = 358
$ws.Range("S6").Value = -461
$ws.Range("T6").Value = 57
$ws.Range("U6").Value = 68
$ws.Range("V6").Value = 1170
$ws.Range("W6").Value = 17.15
$ws.Range("X6").Value = 13.05
$ws.Range("Y6").Value = 7.35
$ws.Range("Z6").Value = 5.01
$ws.Range("AA6").Value = 48.79
$ws.Range("AB6").Value = 5458.18
$ws.Range("AC6").Value = 11532
$ws.Range("AD6").Value = 3.23
$ws.Range("AE6").Value = 171344
$ws.Range("AF6").Value = 0.22
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 1.34
$ws.Range("AI6").Value = 4.17
$ws.Range("AJ6").Value = 2365023

# Remove estimate-year data rows (7,8,9) data beyond column C - kept as error rows
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
